# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the value must be forced to
# text (cells holding numeric-looking strings like "304.96" or "1.000" would
# otherwise be auto-converted to real numbers by Excel).
$updates = @(
    @{ Cell = "D2"; Value = "26.908.91"; ForceText = $false }
    @{ Cell = "E2"; Value = "  -0.76%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "1.862.12"; ForceText = $false }
    @{ Cell = "E3"; Value = "  -0.41%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  -0.09%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "304.96"; ForceText = $true }
    @{ Cell = "E5"; Value = "  -0.82%  "; ForceText = $false }
    @{ Cell = "D7"; Value = "0.5048"; ForceText = $true }
    @{ Cell = "E7"; Value = "  -0.32%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.3623"; ForceText = $true }
    @{ Cell = "E8"; Value = "  -3.37%  "; ForceText = $false }
    @{ Cell = "E9"; Value = "  +0.32%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.8955"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +0.71%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "20.68"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +0.03%  "; ForceText = $false }
    @{ Cell = "D12"; Value = "0.07471"; ForceText = $true }
    @{ Cell = "E12"; Value = "  -1.17%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "1.848.41"; ForceText = $false }
    @{ Cell = "E13"; Value = "  -1.17%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "92.92"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +4.03%  "; ForceText = $false }
    @{ Cell = "D15"; Value = "5.229"; ForceText = $true }
    @{ Cell = "E15"; Value = "  -1.74%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "1.000"; ForceText = $true }
    @{ Cell = "E16"; Value = "  -0.10%  "; ForceText = $false }
    @{ Cell = "D17"; Value = "0.000008471"; ForceText = $true }
    @{ Cell = "E17"; Value = "  -0.25%  "; ForceText = $false }
    @{ Cell = "E18"; Value = "  +0.13%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "1.000"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +0.02%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "26.941.75"; ForceText = $false }
    @{ Cell = "E20"; Value = "  -0.86%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "5.027"; ForceText = $true }
    @{ Cell = "E21"; Value = "  -0.95%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "2.085.83"; ForceText = $false }
    @{ Cell = "E22"; Value = "  -0.84%  "; ForceText = $false }
    @{ Cell = "E23"; Value = "  -2.48%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "6.416"; ForceText = $true }
    @{ Cell = "E24"; Value = "  -1.01%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "147.89"; ForceText = $true }
    @{ Cell = "E25"; Value = "  -1.98%  "; ForceText = $false }
    @{ Cell = "E26"; Value = "  -2.43%  "; ForceText = $false }
    @{ Cell = "D27"; Value = "17.89"; ForceText = $true }
    @{ Cell = "E27"; Value = "  -0.70%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "2.058"; ForceText = $true }
    @{ Cell = "E28"; Value = "  -1.83%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "113.10"; ForceText = $true }
    @{ Cell = "E29"; Value = "  +0.20%  "; ForceText = $false }
    @{ Cell = "D30"; Value = "4.676"; ForceText = $true }
    @{ Cell = "E30"; Value = "  -1.81%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "4.679"; ForceText = $true }
    @{ Cell = "E31"; Value = "  -0.08%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "0.09266"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +2.78%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  -1.01%  "; ForceText = $false }
    @{ Cell = "D34"; Value = "3.002"; ForceText = $true }
    @{ Cell = "E34"; Value = "  -3.05%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "0.7433"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +0.15%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "1.149"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -0.88%  "; ForceText = $false }
    @{ Cell = "D37"; Value = "3.278"; ForceText = $true }
    @{ Cell = "E37"; Value = "  +7.68%  "; ForceText = $false }
    @{ Cell = "B38"; Value = "VeChain"; ForceText = $false }
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; ForceText = $false }
    @{ Cell = "D38"; Value = "0.02002"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -1.73%  "; ForceText = $false }
    @{ Cell = "B39"; Value = "RenderToken"; ForceText = $false }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = $false }
    @{ Cell = "D39"; Value = "2.505"; ForceText = $true }
    @{ Cell = "E39"; Value = "  -1.12%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "0.5579"; ForceText = $true }
    @{ Cell = "E40"; Value = "  +3.57%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "1.071"; ForceText = $true }
    @{ Cell = "E41"; Value = "  -0.45%  "; ForceText = $false }
    @{ Cell = "D42"; Value = "118.25"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +2.57%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "6.475"; ForceText = $true }
    @{ Cell = "E43"; Value = "  -1.67%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "8.525"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +1.16%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "0.1467"; ForceText = $true }
    @{ Cell = "E45"; Value = "  -0.68%  "; ForceText = $false }
    @{ Cell = "D46"; Value = "0.4711"; ForceText = $true }
    @{ Cell = "E46"; Value = "  +1.45%  "; ForceText = $false }
    @{ Cell = "D47"; Value = "1.000"; ForceText = $true }
    @{ Cell = "D48"; Value = "10.03"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +0.44%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "1.562"; ForceText = $true }
    @{ Cell = "E49"; Value = "  -0.35%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "37.04"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +1.29%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "62.97"; ForceText = $true }
    @{ Cell = "E51"; Value = "  -2.50%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Quote-prefix forces Excel to keep the literal text instead of parsing it
        # as a number; resetting the style afterwards drops the quote-prefix format
        # so the cell keeps its original (default) style, same as before the edit.
        $range.Value = "'" + $u.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $u.Value
    }
}
